$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 510

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.1542581303010235
$wsLinear.Range("B3").Value = 0.03394601265597957
$wsLinear.Range("B4").Value = 51.07955614817614
$wsLinear.Range("B5").Value = "[1.0, 0.14524713620753002, 0.026048222905760237, 0.08826587790657203, 0.08683521874205198, 0.020265338434351257, 0.08492595576956288, 0.2723303106873136, 0.10063556677572534, -0.013719874393382533, 0.09823279943222611, 0.0749176217330113, -0.023640790106412977, 0.09231822970644903, 0.2664087422091209, 0.06569849735376455, -0.012008672776967204, 0.11163698280494512, 0.09667318905904564, -0.024479232437878935]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.9497334348819497
$wsNonLinear.Range("B4").Value = 0.009862873733637163
$wsNonLinear.Range("B5").Value = 0.05731836445304229
$wsNonLinear.Range("B6").Value = 48.48729692616877
$wsNonLinear.Range("B7").Value = -0.3595820537345226
$wsNonLinear.Range("B8").Value = 0.08672935989123372
$wsNonLinear.Range("B9").Value = 53.57900730172615
$wsNonLinear.Range("B10").Value = "[1.0, 0.14540705342024962, 0.02564176608172032, 0.08969503440187925, 0.08866458916178428, 0.021576252266940658, 0.0849997147986968, 0.27191586656929645, 0.10069841203064302, -0.01306485732965578, 0.0987281848260582, 0.07569491142622937, -0.022743987051284557, 0.09279042506001473, 0.2660341255975632, 0.06497331774383742, -0.011524208465160474, 0.11245179860214816, 0.0974274160492259, -0.023856607475829135]"
